$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(8).Insert()

$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(8, 3).Value = "Ñuble"
$ws.Cells.Item(8, 4).Value = 44895
$ws.Cells.Item(8, 5).Value = 16
$ws.Cells.Item(8, 6).Value = 100114007
$ws.Cells.Item(8, 7).Value = "Jengibre"
$ws.Cells.Item(8, 8).Value = "Sin especificar"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 30
$ws.Cells.Item(8, 11).Value = 18000
$ws.Cells.Item(8, 12).Value = 18000
$ws.Cells.Item(8, 13).Value = 18000
$ws.Cells.Item(8, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(8, 15).Value = "Perú"
$ws.Cells.Item(8, 16).Value = 1385
$ws.Cells.Item(8, 17).Value = 13
$ws.Cells.Item(8, 18).Value = "Hortaliza"
